$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the cells we touch to be treated as text so Excel does not
# reinterpret numeric-looking strings (e.g. "409.35") as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.741.36"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "3.416.25"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "409.35"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "128.66"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").Value = "43.40"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "0.0000221"
$ws.Range("E12").Value = "  +11.65%  "
$ws.Range("D13").Value = "9.24"
$ws.Range("E13").Value = "  +5.31%  "
$ws.Range("D14").Value = "3.962.00"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "21.10"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").Value = "3.411.76"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "12.34"
$ws.Range("E18").Value = "  +7.79%  "
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "61.760.85"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "481.29"
$ws.Range("E21").Value = "  +28.51%  "
$ws.Range("D22").Value = "91.50"
$ws.Range("E22").Value = "  +4.81%  "
$ws.Range("D23").Value = "3.31"
$ws.Range("E23").Value = "  +3.90%  "
$ws.Range("D24").Value = "13.55"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("D26").Value = "34.48"
$ws.Range("E26").Value = "  +8.59%  "
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  +9.27%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.77"
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.60"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "12.11"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("D31").Value = "0.114"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "41.74"
$ws.Range("E33").Value = "  -5.35%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "59.10"
$ws.Range("E35").Value = "  +13.38%  "
$ws.Range("D36").Value = "0.0499"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "3.46"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").Value = "2.75"
$ws.Range("E39").Value = "  +17.95%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "145.85"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.318"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").Value = "2.08"
$ws.Range("E44").Value = "  +5.10%  "
$ws.Range("D45").Value = "4.34"
$ws.Range("E45").Value = "  +8.64%  "
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  +22.05%  "
$ws.Range("D47").Value = "16.68"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "23.00"
$ws.Range("E48").Value = "  +5.69%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "117.53"
$ws.Range("E49").Value = "  +26.33%  "
$ws.Range("D50").Value = "0.142"
$ws.Range("E50").Value = "  +16.22%  "
$ws.Range("D51").Value = "2.132.60"
$ws.Range("E51").Value = "  +0.83%  "

# Restore the original (default) cell style now that the text values are set.
$ws.Range("B2:E51").Style = "Normal"
